# Rows in this stock report were shuffled: within each small group of
# adjacent line items (same SKU), the Batch No (B), Sale Rate (E),
# Qty (F) and Value (G) columns were redistributed across the rows of
# the group while the SI No (A), Item description (C) and Rate (D)
# stayed put. Apply the new values cell by cell, group by group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B154").Value = 53925
$ws.Range("E154").Value = 79.37
$ws.Range("F154").Value = 1
$ws.Range("G154").Value = 66.44

$ws.Range("B155").Value = 64350
$ws.Range("E155").Value = 70.63
$ws.Range("F155").Value = 101
$ws.Range("G155").Value = 6710.44

$ws.Range("B156").Value = 57756
$ws.Range("F156").Value = -100
$ws.Range("G156").Value = -6644

$ws.Range("B176").Value = 64329
$ws.Range("E176").Value = 128.32
$ws.Range("F176").Value = 6
$ws.Range("G176").Value = 724.14

$ws.Range("B177").Value = 57552
$ws.Range("E177").Value = 136.86
$ws.Range("F177").Value = -5
$ws.Range("G177").Value = -603.45

$ws.Range("B256").Value = 64979
$ws.Range("E256").Value = 314.41
$ws.Range("F256").Value = 82
$ws.Range("G256").Value = 24251.5

$ws.Range("B257").Value = 48719
$ws.Range("E257").Value = 353.35
$ws.Range("F257").Value = -81
$ws.Range("G257").Value = -23955.75

$ws.Range("B271").Value = 48706
$ws.Range("E271").Value = 39.8
$ws.Range("F271").Value = -144
$ws.Range("G271").Value = -4795.2

$ws.Range("B272").Value = 64973
$ws.Range("E272").Value = 35.4
$ws.Range("F272").Value = 150
$ws.Range("G272").Value = 4995

$ws.Range("B309").Value = 63565
$ws.Range("E309").Value = 109.19
$ws.Range("F309").Value = 60
$ws.Range("G309").Value = 6162.6

$ws.Range("B310").Value = 61610
$ws.Range("E310").Value = 122.71
$ws.Range("F310").Value = -58
$ws.Range("G310").Value = -5957.18

$ws.Range("B338").Value = 55373
$ws.Range("E338").Value = 163.62
$ws.Range("F338").Value = -94
$ws.Range("G338").Value = -13562.32

$ws.Range("B339").Value = 63520
$ws.Range("E339").Value = 153.4
$ws.Range("F339").Value = 97
$ws.Range("G339").Value = 13995.16

$ws.Range("B343").Value = 63571
$ws.Range("E343").Value = 152.53
$ws.Range("F343").Value = 29
$ws.Range("G343").Value = 4160.92

$ws.Range("B344").Value = 57802
$ws.Range("E344").Value = 162.71
$ws.Range("F344").Value = -79
$ws.Range("G344").Value = -11334.92

$ws.Range("B367").Value = 63563
$ws.Range("E367").Value = 119.04
$ws.Range("F367").Value = 15
$ws.Range("G367").Value = 1679.4

$ws.Range("B368").Value = 61605
$ws.Range("E368").Value = 133.78
$ws.Range("F368").Value = -13
$ws.Range("G368").Value = -1455.48

$ws.Range("B371").Value = 61608
$ws.Range("E371").Value = 154.12
$ws.Range("F371").Value = -56
$ws.Range("G371").Value = -7224.56

$ws.Range("B372").Value = 63564
$ws.Range("E372").Value = 137.16
$ws.Range("F372").Value = 57
$ws.Range("G372").Value = 7353.57

$ws.Range("B381").Value = 62865
$ws.Range("F381").Value = 151
$ws.Range("G381").Value = 12051.31

$ws.Range("B382").Value = 57817
$ws.Range("F382").Value = 3
$ws.Range("G382").Value = 239.43

$ws.Range("B392").Value = 57835
$ws.Range("F392").Value = 1
$ws.Range("G392").Value = 59.13

$ws.Range("B393").Value = 62933
$ws.Range("F393").Value = 146
$ws.Range("G393").Value = 8632.98

$ws.Range("B411").Value = 63007
$ws.Range("F411").Value = 984
$ws.Range("G411").Value = 168588.72

$ws.Range("B412").Value = 57856
$ws.Range("F412").Value = 2
$ws.Range("G412").Value = 342.66

$ws.Range("B413").Value = 63008
$ws.Range("F413").Value = 504
$ws.Range("G413").Value = 76189.67999999999

$ws.Range("B414").Value = 57857
$ws.Range("F414").Value = 3
$ws.Range("G414").Value = 453.51

$ws.Range("B575").Value = 53263
$ws.Range("E575").Value = 15.29
$ws.Range("F575").Value = -309
$ws.Range("G575").Value = -3958.29

$ws.Range("B576").Value = 65066
$ws.Range("E576").Value = 13.61
$ws.Range("F576").Value = 313
$ws.Range("G576").Value = 4009.53

$ws.Range("B582").Value = 45706
$ws.Range("E582").Value = 23.58
$ws.Range("F582").Value = -202
$ws.Range("G582").Value = -3985.46

$ws.Range("B583").Value = 64922
$ws.Range("E583").Value = 20.98
$ws.Range("F583").Value = 207
$ws.Range("G583").Value = 4084.11

$ws.Range("B585").Value = 64927
$ws.Range("E585").Value = 17.26
$ws.Range("F585").Value = 295
$ws.Range("G585").Value = 4784.9

$ws.Range("B586").Value = 45718
$ws.Range("E586").Value = 19.38
$ws.Range("F586").Value = -294
$ws.Range("G586").Value = -4768.68

$ws.Range("B591").Value = 45709
$ws.Range("E591").Value = 15.69
$ws.Range("F591").Value = -300
$ws.Range("G591").Value = -3945

$ws.Range("B592").Value = 64925
$ws.Range("E592").Value = 13.97
$ws.Range("F592").Value = 302
$ws.Range("G592").Value = 3971.3

$ws.Range("B596").Value = 65067
$ws.Range("E596").Value = 15.65
$ws.Range("F596").Value = 338
$ws.Range("G596").Value = 4978.74

$ws.Range("B597").Value = 53595
$ws.Range("E597").Value = 17.61
$ws.Range("F597").Value = -335
$ws.Range("G597").Value = -4934.55

$ws.Range("B679").Value = 64810
$ws.Range("E679").Value = 291.22
$ws.Range("F679").Value = 7
$ws.Range("G679").Value = 1917.44

$ws.Range("B680").Value = 53319
$ws.Range("E680").Value = 310.64
$ws.Range("F680").Value = -6
$ws.Range("G680").Value = -1643.52
